$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.220588684082031
$ws.Range("B1").Value = 2.022544622421265
$ws.Range("C1").Value = 4.369124412536621
$ws.Range("D1").Value = 2.979888200759888
$ws.Range("E1").Value = 1.188058257102966
